$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.658.34"
$ws.Range("E2").Value = "  -3.26%  "
$ws.Range("D3").Value = "2.088.84"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.79"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5158"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4399"
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09253"
$ws.Range("E9").Value = "  +2.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.06"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.19"
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("D13").Value = "2.084.05"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.754"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.178"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.92"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001157"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.009"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.01"
$ws.Range("E19").Value = "  +8.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06627"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.193"
$ws.Range("D23").Value = "29.720.61"
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.67"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.314"
$ws.Range("E25").Value = "  -3.06%  "
$ws.Range("D26").Value = "2.339.37"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.86"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.19"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.519"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.48"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.141"
$ws.Range("E31").Value = "  -3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1051"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.632"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.183"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.964"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.065"
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.34"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02573"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06719"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.47"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2246"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6834"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6635"
$ws.Range("E44").Value = "  +3.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.14"
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.315"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.615"
$ws.Range("E47").Value = "  -3.74%  "
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000338"
$ws.Range("E49").Value = "  -6.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.64"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("E51").Value = "  -2.15%  "

Write-Host "Applied cryptos list update"
